$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" column header in H1, copying the exact formatting (style)
# used by the other header cells (e.g. G1 "sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Add the numeric value for the new "Save" column in row 2.
$ws.Range("H2").Value = 0
